$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66: 2025-11-02 (serial 45963), 四方坪站 station
$ws.Cells.Item(66, 1).Value = 45963
$ws.Cells.Item(66, 2).Value = "四方坪站"
$ws.Cells.Item(66, 3).Value = 8785.66
$ws.Cells.Item(66, 4).Value = 7812.95
$ws.Cells.Item(66, 5).Value = 2861.71
$ws.Cells.Item(66, 6).Value = 368

# Row 67: 2025-11-02 (serial 45963), 高岭站 station
$ws.Cells.Item(67, 1).Value = 45963
$ws.Cells.Item(67, 2).Value = "高岭站"
$ws.Cells.Item(67, 3).Value = 4392.14
$ws.Cells.Item(67, 4).Value = 3809.79
$ws.Cells.Item(67, 5).Value = 1125.26
$ws.Cells.Item(67, 6).Value = 160

# Update selection to match new last row
$ws.Range("K67").Select()
